$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-26 were re-sorted/re-labelled (4 new "melhores rstds" models moved to
# the top) and every row's metric columns (B:I) were recomputed to the same
# refreshed set of values.
$newNames = @(
    "model_12_3_0",
    "model_12_3_22",
    "model_12_3_21",
    "model_12_3_20",
    "model_12_3_19",
    "model_12_3_18",
    "model_12_3_17",
    "model_12_3_16",
    "model_12_3_15",
    "model_12_3_14",
    "model_12_3_13",
    "model_12_3_23",
    "model_12_3_12",
    "model_12_3_10",
    "model_12_3_9",
    "model_12_3_8",
    "model_12_3_7",
    "model_12_3_6",
    "model_12_3_5",
    "model_12_3_4",
    "model_12_3_3",
    "model_12_3_2",
    "model_12_3_1",
    "model_12_3_11",
    "model_12_3_24"
)

$metrics = @(0.353672031788087, -0.55836533523947, 0.3303799136176995, 0.2441725364439876, 0.7152945399284363, 0.3218813836574554, 0.5454630851745605, 0.4270963668823242)

for ($i = 0; $i -lt $newNames.Count; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $newNames[$i]

    $rowArr = New-Object 'object[,]' 1,8
    for ($j = 0; $j -lt $metrics.Count; $j++) {
        $rowArr[0, $j] = $metrics[$j]
    }
    $ws.Range("B" + $r + ":I" + $r).Value = $rowArr
}
